$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "30.154.09"
Set-TextValue "E2" "  +4.53%  "
Set-TextValue "D3" "1.907.97"
Set-TextValue "E3" "  +5.10%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.28%  "
Set-TextValue "D5" "250.93"
Set-TextValue "E5" "  -0.10%  "
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.21%  "
Set-TextValue "D7" "0.5097"
Set-TextValue "E7" "  +2.37%  "
Set-TextValue "D8" "45.43"
Set-TextValue "E8" "  +4.76%  "
Set-TextValue "D9" "0.2943"
Set-TextValue "E9" "  +5.10%  "
Set-TextValue "D10" "0.06751"
Set-TextValue "E10" "  +5.18%  "
Set-TextValue "D11" "1.913.23"
Set-TextValue "E11" "  +5.19%  "
Set-TextValue "D12" "17.22"
Set-TextValue "E12" "  +2.40%  "
Set-TextValue "D13" "0.07334"
Set-TextValue "E13" "  +2.25%  "
Set-TextValue "D14" "0.6877"
Set-TextValue "E14" "  +5.36%  "
Set-TextValue "D15" "86.48"
Set-TextValue "E15" "  +2.74%  "
Set-TextValue "D16" "4.839"
Set-TextValue "E16" "  +2.19%  "
Set-TextValue "D17" "30.174.80"
Set-TextValue "E17" "  +4.59%  "
Set-TextValue "D18" "0.000008053"
Set-TextValue "E18" "  +8.33%  "
Set-TextValue "D19" "1.002"
Set-TextValue "E19" "  +0.34%  "
Set-TextValue "D20" "12.91"
Set-TextValue "E20" "  +4.79%  "
Set-TextValue "D21" "2.158.92"
Set-TextValue "E21" "  +5.28%  "
Set-TextValue "D22" "1.002"
Set-TextValue "E22" "  +0.27%  "
Set-TextValue "D23" "4.808"
Set-TextValue "E23" "  +3.84%  "
Set-TextValue "D24" "5.697"
Set-TextValue "E24" "  +6.01%  "
Set-TextValue "D25" "9.085"
Set-TextValue "E25" "  +1.74%  "
Set-TextValue "D26" "147.15"
Set-TextValue "E26" "  +2.20%  "
Set-TextValue "D27" "134.35"
Set-TextValue "E27" "  +1.32%  "
Set-TextValue "D28" "17.03"
Set-TextValue "E28" "  +3.75%  "
Set-TextValue "D29" "1.990"
Set-TextValue "E29" "  +4.50%  "
Set-TextValue "D30" "1.397"
Set-TextValue "E30" "  -0.75%  "
Set-TextValue "D31" "4.209"
Set-TextValue "E31" "  +0.62%  "
Set-TextValue "D32" "0.08759"
Set-TextValue "E32" "  +4.27%  "
Set-TextValue "D33" "3.972"
Set-TextValue "E33" "  +2.29%  "
Set-TextValue "D34" "0.05056"
Set-TextValue "E34" "  +1.56%  "
Set-TextValue "D35" "1.144"
Set-TextValue "E35" "  +4.60%  "
Set-TextValue "D36" "0.7099"
Set-TextValue "E36" "  +3.83%  "
Set-TextValue "D37" "2.694"
Set-TextValue "E37" "  +0.70%  "
Set-TextValue "D38" "2.814"
Set-TextValue "E38" "  +0.87%  "
Set-TextValue "D39" "2.279"
Set-TextValue "E39" "  +1.72%  "
Set-TextValue "D40" "0.9678"
Set-TextValue "E40" "  -0.13%  "
Set-TextValue "D41" "0.01678"
Set-TextValue "E41" "  +5.11%  "
Set-TextValue "D42" "6.045"
Set-TextValue "E42" "  -0.64%  "
Set-TextValue "D43" "104.99"
Set-TextValue "E43" "  +3.95%  "
Set-TextValue "D44" "0.4258"
Set-TextValue "E44" "  +3.08%  "
Set-TextValue "E45" "  -0.04%  "
Set-TextValue "D46" "7.568"
Set-TextValue "E46" "  +4.02%  "
Set-TextValue "D47" "0.1269"
Set-TextValue "E47" "  +3.18%  "
Set-TextValue "D48" "0.05743"
Set-TextValue "E48" "  +3.94%  "
Set-TextValue "B49" "Elrond"
Set-TextValue "C49" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D49" "32.96"
Set-TextValue "E49" "  +3.63%  "
Set-TextValue "B50" "EnergySwap"
Set-TextValue "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "8.430"
Set-TextValue "E50" "  +2.74%  "
Set-TextValue "D51" "0.3786"
Set-TextValue "E51" "  +3.51%  "
